# Generate Report for Handback
# This mirrors the refresh of the handback status report: the first
# handback batch's identifiers/hashes/timestamps are replaced with a new
# batch, and the second file's timestamps now collapse onto the same
# (new) handoff/handback record as the first file.

$wb = $excel.ActiveWorkbook

$oldMdUuid1 = "2bc49e7c-1f69-4273-ba0d-714a75eb899d.md"
$newMdUuid1 = "0646e59e-8ec3-4821-a9de-21864af14c32.md"

$oldMdUuid2 = "e7faba28-3992-4225-bdb9-c5f7d617bd3e.md"
$newMdUuid2 = "ffff17bc4d83-0e98-4527-855f-5854249b09e1.md"

$newZhXlf = "0646e59e-8ec3-4821-a9de-21864af14c32.d38632e3469738437b3153b9189d010b22a7957f.zh-cn.xlf"
$newDeXlf = "0646e59e-8ec3-4821-a9de-21864af14c32.d38632e3469738437b3153b9189d010b22a7957f.de-de.xlf"

$newZhHandoffDt = "2016-03-20 00:49:14"
$newZhHandbackDt = "2016-03-20 00:49:31"
$newDeHandoffDt = "2016-03-20 00:49:17"
$newDeHandbackDt = "2016-03-20 00:49:37"

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newMdUuid1
$wsOverview.Range("A3").Value = $newMdUuid2

foreach ($hl in $wsOverview.Hyperlinks) {
    if ($hl.TextToDisplay -eq $oldMdUuid1) {
        $hl.TextToDisplay = $newMdUuid1
    } elseif ($hl.TextToDisplay -eq $oldMdUuid2) {
        $hl.TextToDisplay = $newMdUuid2
    }
}

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = $newMdUuid1
$wsZh.Range("F2").Value = $newMdUuid1
$wsZh.Range("D2").Value = $newZhXlf
$wsZh.Range("G2").Value = $newZhXlf
$wsZh.Range("E2").Value = $newZhHandoffDt
$wsZh.Range("H2").Value = $newZhHandbackDt

$wsZh.Range("A3").Value = $newMdUuid2
$wsZh.Range("F3").Value = $newMdUuid2
$wsZh.Range("D3").Value = $newZhXlf
$wsZh.Range("G3").Value = $newZhXlf
$wsZh.Range("E3").Value = $newZhHandoffDt
$wsZh.Range("H3").Value = $newZhHandbackDt

foreach ($hl in $wsZh.Hyperlinks) {
    if ($hl.TextToDisplay -eq $oldMdUuid1) {
        $hl.TextToDisplay = $newMdUuid1
    } elseif ($hl.TextToDisplay -eq $oldMdUuid2) {
        $hl.TextToDisplay = $newMdUuid2
    } elseif ($hl.TextToDisplay -eq "2bc49e7c-1f69-4273-ba0d-714a75eb899d.8a4e2a08faba8868e51f71545f50f9b9092f20f5.zh-cn.xlf") {
        $hl.TextToDisplay = $newZhXlf
    } elseif ($hl.TextToDisplay -eq "e7faba28-3992-4225-bdb9-c5f7d617bd3e.22b7535328d6b6d3b62c784bcce841610a6765dd.zh-cn.xlf") {
        $hl.TextToDisplay = $newZhXlf
    }
}

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = $newMdUuid1
$wsDe.Range("F2").Value = $newMdUuid1
$wsDe.Range("D2").Value = $newDeXlf
$wsDe.Range("G2").Value = $newDeXlf
$wsDe.Range("E2").Value = $newDeHandoffDt
$wsDe.Range("H2").Value = $newDeHandbackDt

$wsDe.Range("A3").Value = $newMdUuid2
$wsDe.Range("F3").Value = $newMdUuid2
$wsDe.Range("D3").Value = $newDeXlf
$wsDe.Range("G3").Value = $newDeXlf
$wsDe.Range("E3").Value = $newDeHandoffDt
$wsDe.Range("H3").Value = $newDeHandbackDt

foreach ($hl in $wsDe.Hyperlinks) {
    if ($hl.TextToDisplay -eq $oldMdUuid1) {
        $hl.TextToDisplay = $newMdUuid1
    } elseif ($hl.TextToDisplay -eq $oldMdUuid2) {
        $hl.TextToDisplay = $newMdUuid2
    } elseif ($hl.TextToDisplay -eq "2bc49e7c-1f69-4273-ba0d-714a75eb899d.8a4e2a08faba8868e51f71545f50f9b9092f20f5.de-de.xlf") {
        $hl.TextToDisplay = $newDeXlf
    } elseif ($hl.TextToDisplay -eq "e7faba28-3992-4225-bdb9-c5f7d617bd3e.22b7535328d6b6d3b62c784bcce841610a6765dd.de-de.xlf") {
        $hl.TextToDisplay = $newDeXlf
    }
}
